# Applies stock-report quantity/value corrections described by the commit diff.
# All edits are plain value updates (no formulas exist in this worksheet);
# Column F = Qty, Column G = Value (Rate * Qty), Column B on "Sub Total"/"Grand
# Total" rows = rollup of the Value column for that company block / whole sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F74").Value = 14
$ws.Range("G74").Value = 984.48
$ws.Range("F80").Value = 21
$ws.Range("G80").Value = 1305.57
$ws.Range("B83").Value = 107678.32
$ws.Range("F197").Value = 1000
$ws.Range("G197").Value = 18500
$ws.Range("B204").Value = 29136.36
$ws.Range("F222").Value = 19
$ws.Range("G222").Value = 11930.48
$ws.Range("B234").Value = 162186.01
$ws.Range("F247").Value = 148
$ws.Range("G247").Value = 6845
$ws.Range("F256").Value = 39
$ws.Range("G256").Value = 1971.84
$ws.Range("F272").Value = 64
$ws.Range("G272").Value = 5258.24
$ws.Range("F273").Value = 31
$ws.Range("G273").Value = 1566.74
$ws.Range("B276").Value = 266652.03
$ws.Range("F279").Value = 146
$ws.Range("G279").Value = 25266.76
$ws.Range("F283").Value = 74
$ws.Range("G283").Value = 6137.56
$ws.Range("F293").Value = 61
$ws.Range("G293").Value = 9814.9
$ws.Range("F304").Value = 220
$ws.Range("G304").Value = 10076
$ws.Range("F306").Value = 103
$ws.Range("G306").Value = 11622.52
$ws.Range("F310").Value = 111
$ws.Range("G310").Value = 9221.879999999999
$ws.Range("F313").Value = 127
$ws.Range("G313").Value = 7489.19
$ws.Range("F314").Value = 324
$ws.Range("G314").Value = 4393.44
$ws.Range("F316").Value = 153
$ws.Range("G316").Value = 2238.39
$ws.Range("F331").Value = 87
$ws.Range("G331").Value = 10811.49
$ws.Range("F337").Value = 95
$ws.Range("G337").Value = 8143.4
$ws.Range("F348").Value = 15
$ws.Range("G348").Value = 1170.75
$ws.Range("F354").Value = 42
$ws.Range("G354").Value = 4364.22
$ws.Range("F355").Value = 1
$ws.Range("G355").Value = 118.63
$ws.Range("F362").Value = 944
$ws.Range("G362").Value = 161735.52
$ws.Range("F363").Value = 450
$ws.Range("G363").Value = 68026.5
$ws.Range("F364").Value = 60
$ws.Range("G364").Value = 25109.4
$ws.Range("F372").Value = 0
$ws.Range("G372").Value = 0
$ws.Range("F374").Value = 758
$ws.Range("G374").Value = 45078.26
$ws.Range("B380").Value = 1249854.77
$ws.Range("F418").Value = 39
$ws.Range("G418").Value = 5855.46
$ws.Range("B422").Value = 48447.58
$ws.Range("F467").Value = 4
$ws.Range("G467").Value = 241.8
$ws.Range("B491").Value = 64831.98
$ws.Range("F503").Value = 14
$ws.Range("G503").Value = 39907
$ws.Range("B518").Value = 333261.17
$ws.Range("F569").Value = 54
$ws.Range("G569").Value = 4404.24
$ws.Range("B570").Value = 42946
$ws.Range("F583").Value = 8
$ws.Range("G583").Value = 1421.28
$ws.Range("B584").Value = 1489.81
$ws.Range("F608").Value = 212
$ws.Range("G608").Value = 3362.32
$ws.Range("B609").Value = 3991.41
$ws.Range("F663").Value = 117
$ws.Range("G663").Value = 8372.52
$ws.Range("F665").Value = 178
$ws.Range("G665").Value = 16331.5
$ws.Range("F667").Value = 103
$ws.Range("G667").Value = 6275.79
$ws.Range("F668").Value = 97
$ws.Range("G668").Value = 5910.21
$ws.Range("F669").Value = 107
$ws.Range("G669").Value = 7656.92
$ws.Range("F670").Value = 44
$ws.Range("G670").Value = 6287.16
$ws.Range("F672").Value = 139
$ws.Range("G672").Value = 9946.84
$ws.Range("B673").Value = 99516
$ws.Range("F709").Value = 512
$ws.Range("G709").Value = 41758.72
$ws.Range("F711").Value = 25
$ws.Range("G711").Value = 2039
$ws.Range("F712").Value = 264
$ws.Range("G712").Value = 40793.28
$ws.Range("F713").Value = 287
$ws.Range("G713").Value = 23407.72
$ws.Range("F718").Value = 318
$ws.Range("G718").Value = 6906.96
$ws.Range("F720").Value = 35
$ws.Range("G720").Value = 3173.8
$ws.Range("F721").Value = 422
$ws.Range("G721").Value = 29354.32
$ws.Range("F729").Value = 956
$ws.Range("G729").Value = 137664
$ws.Range("F730").Value = 23
$ws.Range("G730").Value = 860.66
$ws.Range("F731").Value = 766
$ws.Range("G731").Value = 92463.86
$ws.Range("B733").Value = 703792.97
$ws.Range("F773").Value = 72
$ws.Range("G773").Value = 2176.56
$ws.Range("F775").Value = 192
$ws.Range("G775").Value = 8768.639999999999
$ws.Range("F778").Value = 2977
$ws.Range("G778").Value = 485578.47
$ws.Range("F780").Value = 1
$ws.Range("G780").Value = 144.65
$ws.Range("F782").Value = 5
$ws.Range("G782").Value = 642.9
$ws.Range("B783").Value = 508308.46
$ws.Range("B795").Value = 6266085.86
$ws.Range("B796").Value = 6266085.86
